$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -3
$ws.Range("F8").Value = -9
$ws.Range("F10").Value = -10
$ws.Range("F11").Value = -4
$ws.Range("F14").Value = -4
$ws.Range("F15").Value = -2
$ws.Range("F17").Value = -4
$ws.Range("F18").Value = -6
$ws.Range("F22").Value = -3
